# Regenerate instances - apply updated randomly-generated values produced
# by the create script for the 01_Uniform instance workbook.

$wb = $excel.ActiveWorkbook

# ----- Productdata sheet -----
$ws = $wb.Worksheets.Item("Productdata")

$ws.Range("E2").Value = 0.009750000000000002
$ws.Range("E3").Value = 0.009300000000000001

$ws.Range("C4").Value = 4
$ws.Range("E4").Value = 0.0036

$ws.Range("C5").Value = 4
$ws.Range("E5").Value = 0.0015

$ws.Range("C6").Value = 4
$ws.Range("E6").Value = 0.0027

$ws.Range("C7").Value = 1
$ws.Range("E7").Value = 0.004875000000000001

$ws.Range("C8").Value = 1
$ws.Range("E8").Value = 0.009525

$ws.Range("C9").Value = 1
$ws.Range("E9").Value = 0.00465

# ----- Capacity sheet -----
$ws = $wb.Worksheets.Item("Capacity")

$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 20
$ws.Range("B6").Value = 40
$ws.Range("B7").Value = 25
$ws.Range("B9").Value = 20

# ----- ProcessingTime sheet -----
$ws = $wb.Worksheets.Item("ProcessingTime")

$ws.Range("B2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D4").Value = 1
$ws.Range("F6").Value = 2
$ws.Range("G7").Value = 5
$ws.Range("I9").Value = 4
